$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates (Sostdc1 -> Lrp6, MuSCs -> MuSCs)
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 12.78754244675556
$ws.Range("R2").Value = 115.0878820208
$ws.Range("S2").Value = 0.1720325859617629
$ws.Range("T2").Value = 0.1720325859617629

# Row 3 updates (Sostdc1 -> Lrp6, MuSCs -> ECs)
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("Q3").Value = 45.32133362026111
$ws.Range("R3").Value = 407.89200258235
$ws.Range("S3").Value = 0.6097142007069145
$ws.Range("T3").Value = 0.6097142007069145

# Row 4 updates (Sostdc1 -> Lrp6, MuSCs -> MuSCs)
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("S4").Value = 0.2182532133313226
$ws.Range("T4").Value = 0.2182532133313226
